# Weekly CompStat report refresh: new crime data collected (week of 6/12/2023-6/18/2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: volume/issue number and reporting week dates ---
$ws.Range("A8").Value = "Volume 30   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# --- Update weekly crime statistics table (rows 14-29) ---
$ws.Range("F14").Value = 1
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0"
$ws.Range("C28").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "***.*"
$ws.Range("N26").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("L15").Value = 20
$ws.Range("N15").Value = -60
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 19
$ws.Range("H16").Value = 11.764705882352
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = -6.666666666666
$ws.Range("L16").Value = 62.318840579710
$ws.Range("M16").Value = 6.666666666666
$ws.Range("N16").Value = -69.647696476964
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -5.882352941176
$ws.Range("I17").Value = 170
$ws.Range("J17").Value = 208
$ws.Range("K17").Value = -18.269230769230
$ws.Range("L17").Value = 17.241379310344
$ws.Range("M17").Value = 50.442477876106
$ws.Range("N17").Value = -39.501779359430
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value2 = "0"
$ws.Range("C28").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = -73.333333333333
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = -45.744680851063
$ws.Range("L18").Value = 37.837837837837
$ws.Range("M18").Value = 4.081632653061
$ws.Range("N18").Value = -78.481012658227
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -22.448979591836
$ws.Range("I19").Value = 197
$ws.Range("J19").Value = 231
$ws.Range("K19").Value = -14.718614718614
$ws.Range("L19").Value = 48.120300751879
$ws.Range("M19").Value = 97
$ws.Range("N19").Value = 21.604938271604
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = -18.918918918918
$ws.Range("L20").Value = 11.111111111111
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -80
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -29.411764705882
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = -22.962962962963
$ws.Range("I21").Value = 570
$ws.Range("J21").Value = 699
$ws.Range("K21").Value = -18.454935622317
$ws.Range("L21").Value = 34.751773049645
$ws.Range("M21").Value = 41.089108910891
$ws.Range("N21").Value = -53.620829943043
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value2 = "0"
$ws.Range("C28").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("D28").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("E28").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = -38.888888888888
$ws.Range("L22").Value = -15.384615384615
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value2 = "0"
$ws.Range("C28").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = -65
$ws.Range("J23").Value = 82
$ws.Range("K23").Value = -14.634146341463
$ws.Range("L23").Value = -2.777777777777
$ws.Range("M23").Value = 55.555555555555
$ws.Range("C24").Value = 25
$ws.Range("E24").Value = -10.714285714285
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = -9.401709401709
$ws.Range("I24").Value = 583
$ws.Range("J24").Value = 537
$ws.Range("K24").Value = 8.566108007448
$ws.Range("L24").Value = 77.203647416413
$ws.Range("M24").Value = 41.162227602905
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = 13.953488372093
$ws.Range("I25").Value = 238
$ws.Range("J25").Value = 260
$ws.Range("K25").Value = -8.461538461538
$ws.Range("L25").Value = 21.428571428571
$ws.Range("M25").Value = 10.697674418604
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0"
$ws.Range("C28").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "***.*"
$ws.Range("N26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("L26").Value = 10
$ws.Range("C27").Value = 1
$ws.Range("D28").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 25
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -19.354838709677
$ws.Range("L27").Value = -30.555555555555
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 6.25
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 15
$ws.Range("K29").Value = 6.666666666666

$excel.CutCopyMode = 0
